$d = $word.ActiveDocument

# Phase 1: replace each old value with a unique placeholder to avoid collisions
# (some new values equal old values used elsewhere in the document)
$d.Content.Find.Execute("2024-01-19 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_0@@", 2)
$d.Content.Find.Execute("66÷7=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_1@@", 2)
$d.Content.Find.Execute("64÷6=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_2@@", 2)
$d.Content.Find.Execute("99÷8=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_3@@", 2)
$d.Content.Find.Execute("19÷4=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_4@@", 2)
$d.Content.Find.Execute("58÷3=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_5@@", 2)
$d.Content.Find.Execute("64÷3=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_6@@", 2)
$d.Content.Find.Execute("29÷7=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_7@@", 2)
$d.Content.Find.Execute("45÷8=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_8@@", 2)
$d.Content.Find.Execute("40÷3=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_9@@", 2)
$d.Content.Find.Execute("17÷4=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_10@@", 2)
$d.Content.Find.Execute("96÷9=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_11@@", 2)
$d.Content.Find.Execute("54÷2=27, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_12@@", 2)
$d.Content.Find.Execute("15÷5=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_13@@", 2)
$d.Content.Find.Execute("54÷9=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_14@@", 2)
$d.Content.Find.Execute("82÷3=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_15@@", 2)
$d.Content.Find.Execute("17÷6=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_16@@", 2)
$d.Content.Find.Execute("54÷3=18, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_17@@", 2)
$d.Content.Find.Execute("66÷2=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_18@@", 2)
$d.Content.Find.Execute("70÷2=35, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_19@@", 2)
$d.Content.Find.Execute("25÷9=2, 7", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_20@@", 2)
$d.Content.Find.Execute("68÷4=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_21@@", 2)
$d.Content.Find.Execute("18÷6=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_22@@", 2)
$d.Content.Find.Execute("94÷3=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_23@@", 2)
$d.Content.Find.Execute("53÷3=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_24@@", 2)
$d.Content.Find.Execute("17÷5=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_25@@", 2)

# Phase 2: replace placeholders with final values
$d.Content.Find.Execute("@@PLACEHOLDER_0@@", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-20 Saturday", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_1@@", $true, $false, $false, $false, $false, $true, 1, $false, "53÷3=17, 2", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_2@@", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=11, 7", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_3@@", $true, $false, $false, $false, $false, $true, 1, $false, "27÷2=13, 1", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_4@@", $true, $false, $false, $false, $false, $true, 1, $false, "10÷3=3, 1", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_5@@", $true, $false, $false, $false, $false, $true, 1, $false, "93÷8=11, 5", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_6@@", $true, $false, $false, $false, $false, $true, 1, $false, "92÷8=11, 4", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_7@@", $true, $false, $false, $false, $false, $true, 1, $false, "35÷5=7, 0", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_8@@", $true, $false, $false, $false, $false, $true, 1, $false, "34÷2=17, 0", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_9@@", $true, $false, $false, $false, $false, $true, 1, $false, "93÷9=10, 3", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_10@@", $true, $false, $false, $false, $false, $true, 1, $false, "37÷5=7, 2", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_11@@", $true, $false, $false, $false, $false, $true, 1, $false, "39÷6=6, 3", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_12@@", $true, $false, $false, $false, $false, $true, 1, $false, "56÷9=6, 2", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_13@@", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=9, 1", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_14@@", $true, $false, $false, $false, $false, $true, 1, $false, "91÷7=13, 0", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_15@@", $true, $false, $false, $false, $false, $true, 1, $false, "42÷5=8, 2", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_16@@", $true, $false, $false, $false, $false, $true, 1, $false, "20÷4=5, 0", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_17@@", $true, $false, $false, $false, $false, $true, 1, $false, "87÷4=21, 3", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_18@@", $true, $false, $false, $false, $false, $true, 1, $false, "65÷6=10, 5", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_19@@", $true, $false, $false, $false, $false, $true, 1, $false, "98÷3=32, 2", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_20@@", $true, $false, $false, $false, $false, $true, 1, $false, "83÷9=9, 2", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_21@@", $true, $false, $false, $false, $false, $true, 1, $false, "20÷3=6, 2", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_22@@", $true, $false, $false, $false, $false, $true, 1, $false, "94÷6=15, 4", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_23@@", $true, $false, $false, $false, $false, $true, 1, $false, "87÷4=21, 3", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_24@@", $true, $false, $false, $false, $false, $true, 1, $false, "80÷7=11, 3", 2)
$d.Content.Find.Execute("@@PLACEHOLDER_25@@", $true, $false, $false, $false, $false, $true, 1, $false, "42÷5=8, 2", 2)
